$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix P193: style 13 (unused/odd font) -> style 1 (matches neighboring cells)
$ws.Range("O193").Copy()
$ws.Range("P193").PasteSpecial(-4122)

# Append rows 204-214 with data, cloning row formatting from the matching
# existing pattern rows (195 = "all s=11" rows, 194 = "all s=1" rows)

# Row 204
$ws.Range("A195:Q195").Copy()
$ws.Range("A204:Q204").PasteSpecial(-4122)
$ws.Range("A204").Value = 44154.0
$ws.Range("B204").Value = 958953.0
$ws.Range("C204").Value = 141422.0
$ws.Range("D204").Value = 816833.0
$ws.Range("E204").Value = 698.0
$ws.Range("F204").Value = 2280.0
$ws.Range("G204").Value = 135701.0
$ws.Range("H204").Value = 200.0
$ws.Range("I204").Value = 29.0
$ws.Range("J204").Value = 44.0
$ws.Range("K204").Value = 0.0
$ws.Range("L204").Value = 3.0
$ws.Range("M204").Value = 868.0
$ws.Range("N204").Value = 1208.0
$ws.Range("O204").Value = 39.0
$ws.Range("P204").Value = 123.0
$ws.Range("Q204").Value = 42.0

# Row 205
$ws.Range("A195:Q195").Copy()
$ws.Range("A205:Q205").PasteSpecial(-4122)
$ws.Range("A205").Value = 44155.0
$ws.Range("B205").Value = 962881.0
$ws.Range("C205").Value = 141582.0
$ws.Range("D205").Value = 820082.0
$ws.Range("E205").Value = 1207.0
$ws.Range("F205").Value = 2282.0
$ws.Range("G205").Value = 135748.0
$ws.Range("H205").Value = 188.0
$ws.Range("I205").Value = 25.0
$ws.Range("J205").Value = 45.0
$ws.Range("K205").Value = 0.0
$ws.Range("L205").Value = 3.0
$ws.Range("M205").Value = 870.0
$ws.Range("N205").Value = 1208.0
$ws.Range("O205").Value = 39.0
$ws.Range("P205").Value = 123.0
$ws.Range("Q205").Value = 42.0

# Row 206
$ws.Range("A194:Q194").Copy()
$ws.Range("A206:Q206").PasteSpecial(-4122)
$ws.Range("A206").Value = 44156.0
$ws.Range("B206").Value = 966284.0
$ws.Range("C206").Value = 141758.0
$ws.Range("D206").Value = 823880.0
$ws.Range("E206").Value = 646.0
$ws.Range("F206").Value = 2284.0
$ws.Range("G206").Value = 135999.0
$ws.Range("H206").Value = 189.0
$ws.Range("I206").Value = 21.0
$ws.Range("J206").Value = 44.0
$ws.Range("K206").Value = 0.0
$ws.Range("L206").Value = 3.0
$ws.Range("M206").Value = 871.0
$ws.Range("N206").Value = 1209.0
$ws.Range("O206").Value = 39.0
$ws.Range("P206").Value = 123.0
$ws.Range("Q206").Value = 42.0

# Row 207
$ws.Range("A194:Q194").Copy()
$ws.Range("A207:Q207").PasteSpecial(-4122)
$ws.Range("A207").Value = 44157.0
$ws.Range("B207").Value = 971251.0
$ws.Range("C207").Value = 141917.0
$ws.Range("D207").Value = 828072.0
$ws.Range("E207").Value = 1262.0
$ws.Range("F207").Value = 2284.0
$ws.Range("G207").Value = 136289.0
$ws.Range("H207").Value = 189.0
$ws.Range("I207").Value = 24.0
$ws.Range("J207").Value = 44.0
$ws.Range("K207").Value = 0.0
$ws.Range("L207").Value = 3.0
$ws.Range("M207").Value = 871.0
$ws.Range("N207").Value = 1209.0
$ws.Range("O207").Value = 39.0
$ws.Range("P207").Value = 123.0
$ws.Range("Q207").Value = 42.0

# Row 208
$ws.Range("A194:Q194").Copy()
$ws.Range("A208:Q208").PasteSpecial(-4122)
$ws.Range("A208").Value = 44158.0
$ws.Range("B208").Value = 972921.0
$ws.Range("C208").Value = 141979.0
$ws.Range("D208").Value = 829603.0
$ws.Range("E208").Value = 1339.0
$ws.Range("F208").Value = 2284.0
$ws.Range("G208").Value = 136548.0
$ws.Range("H208").Value = 194.0
$ws.Range("I208").Value = 24.0
$ws.Range("J208").Value = 43.0
$ws.Range("K208").Value = 0.0
$ws.Range("L208").Value = 3.0
$ws.Range("M208").Value = 871.0
$ws.Range("N208").Value = 1209.0
$ws.Range("O208").Value = 39.0
$ws.Range("P208").Value = 123.0
$ws.Range("Q208").Value = 42.0

# Row 209
$ws.Range("A195:Q195").Copy()
$ws.Range("A209:Q209").PasteSpecial(-4122)
$ws.Range("A209").Value = 44159.0
$ws.Range("B209").Value = 977488.0
$ws.Range("C209").Value = 142168.0
$ws.Range("D209").Value = 833268.0
$ws.Range("E209").Value = 2052.0
$ws.Range("F209").Value = 2287.0
$ws.Range("G209").Value = 136777.0
$ws.Range("H209").Value = 186.0
$ws.Range("I209").Value = 24.0
$ws.Range("J209").Value = 42.0
$ws.Range("K209").Value = 0.0
$ws.Range("L209").Value = 5.0
$ws.Range("M209").Value = 871.0
$ws.Range("N209").Value = 1212.0
$ws.Range("O209").Value = 39.0
$ws.Range("P209").Value = 123.0
$ws.Range("Q209").Value = 42.0

# Row 210
$ws.Range("A195:Q195").Copy()
$ws.Range("A210:Q210").PasteSpecial(-4122)
$ws.Range("A210").Value = 44160.0
$ws.Range("B210").Value = 980902.0
$ws.Range("C210").Value = 142379.0
$ws.Range("D210").Value = 838114.0
$ws.Range("E210").Value = 409.0
$ws.Range("F210").Value = 2289.0
$ws.Range("G210").Value = 137007.0
$ws.Range("H210").Value = 175.0
$ws.Range("I210").Value = 20.0
$ws.Range("J210").Value = 40.0
$ws.Range("K210").Value = 0.0
$ws.Range("L210").Value = 5.0
$ws.Range("M210").Value = 871.0
$ws.Range("N210").Value = 1214.0
$ws.Range("O210").Value = 39.0
$ws.Range("P210").Value = 123.0
$ws.Range("Q210").Value = 42.0

# Row 211
$ws.Range("A195:Q195").Copy()
$ws.Range("A211:Q211").PasteSpecial(-4122)
$ws.Range("A211").Value = 44161.0
$ws.Range("B211").Value = 985284.0
$ws.Range("C211").Value = 142623.0
$ws.Range("D211").Value = 841926.0
$ws.Range("E211").Value = 735.0
$ws.Range("F211").Value = 2291.0
$ws.Range("G211").Value = 137175.0
$ws.Range("H211").Value = 178.0
$ws.Range("I211").Value = 20.0
$ws.Range("J211").Value = 40.0
$ws.Range("K211").Value = 0.0
$ws.Range("L211").Value = 4.0
$ws.Range("M211").Value = 873.0
$ws.Range("N211").Value = 1214.0
$ws.Range("O211").Value = 39.0
$ws.Range("P211").Value = 123.0
$ws.Range("Q211").Value = 42.0

# Row 212
$ws.Range("A195:Q195").Copy()
$ws.Range("A212:Q212").PasteSpecial(-4122)
$ws.Range("A212").Value = 44162.0
$ws.Range("B212").Value = 989029.0
$ws.Range("C212").Value = 142936.0
$ws.Range("D212").Value = 845413.0
$ws.Range("E212").Value = 680.0
$ws.Range("F212").Value = 2293.0
$ws.Range("G212").Value = 137229.0
$ws.Range("H212").Value = 183.0
$ws.Range("I212").Value = 19.0
$ws.Range("J212").Value = 42.0
$ws.Range("K212").Value = 0.0
$ws.Range("L212").Value = 3.0
$ws.Range("M212").Value = 873.0
$ws.Range("N212").Value = 1216.0
$ws.Range("O212").Value = 39.0
$ws.Range("P212").Value = 123.0
$ws.Range("Q212").Value = 42.0

# Row 213
$ws.Range("A194:Q194").Copy()
$ws.Range("A213:Q213").PasteSpecial(-4122)
$ws.Range("A213").Value = 44163.0
$ws.Range("B213").Value = 993860.0
$ws.Range("C213").Value = 143393.0
$ws.Range("D213").Value = 848998.0
$ws.Range("E213").Value = 1469.0
$ws.Range("F213").Value = 2294.0
$ws.Range("G213").Value = 137500.0
$ws.Range("H213").Value = 183.0
$ws.Range("I213").Value = 19.0
$ws.Range("J213").Value = 42.0
$ws.Range("K213").Value = 0.0
$ws.Range("L213").Value = 3.0
$ws.Range("M213").Value = 874.0
$ws.Range("N213").Value = 1216.0
$ws.Range("O213").Value = 39.0
$ws.Range("P213").Value = 123.0
$ws.Range("Q213").Value = 42.0

# Row 214
$ws.Range("A194:Q194").Copy()
$ws.Range("A214:Q214").PasteSpecial(-4122)
$ws.Range("A214").Value = 44164.0
$ws.Range("B214").Value = 998249.0
$ws.Range("C214").Value = 143642.0
$ws.Range("D214").Value = 852992.0
$ws.Range("E214").Value = 1615.0
$ws.Range("F214").Value = 2295.0
$ws.Range("G214").Value = 137722.0
$ws.Range("H214").Value = 197.0
$ws.Range("I214").Value = 26.0
$ws.Range("J214").Value = 42.0
$ws.Range("K214").Value = 0.0
$ws.Range("L214").Value = 5.0
$ws.Range("M214").Value = 874.0
$ws.Range("N214").Value = 1217.0
$ws.Range("O214").Value = 39.0
$ws.Range("P214").Value = 123.0
$ws.Range("Q214").Value = 42.0

$excel.CutCopyMode = $false
